$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K) values for rows 2-6 per regenerated save_data
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
